$d = $word.ActiveDocument

# The phrase "React, " appears twice in the resume (once in the
# Professional Summary, once in the bullet list under the Sotreq
# entry). We need to insert "Bootstrap, " only in the bullet listing
# the tech stack ("... React, Node, Python, SQL, MySQL, Qt/QML."), so
# it is located via the unique surrounding text "React, Node".

$findRng = $d.Content
$found = $findRng.Find.Execute(
    "React, Node",   # FindText
    $true,           # MatchCase
    $false,          # MatchWholeWord
    $false,          # MatchWildcards
    $false,          # MatchSoundsLike
    $false,          # MatchAllWordForms
    $true,           # Forward
    1,               # Wrap (wdFindContinue)
    $false,          # Format
    "",              # ReplaceWith
    0                # Replace (wdReplaceNone)
)

if (-not $found) {
    throw "Could not find target text 'React, Node' in the document."
}

$matchStart = $findRng.Start

# Range covering just "React, " (7 characters), immediately before "Node".
$insertionPoint = $d.Range($matchStart, $matchStart + 7)
$insertionPoint.Collapse(0)

# Insert "Bootstrap, " as a tracked-change revision and accept just
# that revision. This mirrors how Word keeps a freshly typed insertion
# in its own run (distinct from the identically-formatted text around
# it) instead of silently re-merging it into the neighboring run, while
# avoiding the side effects of a document-wide "accept all revisions".
$wasTracking = $d.TrackRevisions
$d.TrackRevisions = $true
$insertionPoint.InsertAfter("Bootstrap, ")
$d.TrackRevisions = $wasTracking

if ($d.Revisions.Count -gt 0) {
    $d.Revisions.Item(1).Accept()
}
